{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Change 1: \"...specific technologies like Airflow and Google Cloud\n// Platform.\" -> \"...specific technologies like Apache Airflow and\n// Google Cloud Platform.\"  (rendered as three runs: \"...like \",\n// \"Apache \", \"Airflow and Google Cloud Platform.\")\n// ---------------------------------------------------------------------\nconst airflowResults = body.search(\"Airflow and Google Cloud Platform.\", {\n  matchCase: true,\n});\nairflowResults.load(\"items\");\nawait context.sync();\n\nif (airflowResults.items.length > 0) {\n  // Insert the new word right before \"Airflow ...\".\n  airflowResults.items[0].insertText(\"Apache \", Word.InsertLocation.before);\n  await context.sync();\n\n  // Isolate the just-inserted \"Apache \" word as its own range (between the\n  // end of the unique anchor \"...specific technologies like \" and the\n  // start of the now-shifted \"Airflow and Google Cloud Platform.\") and\n  // touch its font so the engine materialises it as a dedicated run,\n  // matching the run layout produced by the original edit.\n  const beforeAnchor = body.search(\"specific technologies like \", {\n    matchCase: true,\n  });\n  beforeAnchor.load(\"items\");\n  await context.sync();\n\n  if (beforeAnchor.items.length > 0) {\n    const startPoint = beforeAnchor.items[0].getRange(\"End\");\n\n    const airflowResults2 = body.search(\"Airflow and Google Cloud Platform.\", {\n      matchCase: true,\n    });\n    airflowResults2.load(\"items\");\n    await context.sync();\n\n    const endPoint = airflowResults2.items[0].getRange(\"Start\");\n\n    const apacheRange = startPoint.expandTo(endPoint);\n    apacheRange.font.set({ name: \"Aparajita\" });\n    await context.sync();\n\n    // Re-touch the trailing \"Airflow and Google Cloud Platform.\" range too\n    // so it becomes its own run, separate from \"Apache \".\n    const airflowResults3 = body.search(\"Airflow and Google Cloud Platform.\", {\n      matchCase: true,\n    });\n    airflowResults3.load(\"items\");\n    await context.sync();\n    airflowResults3.items[0].font.set({ name: \"Aparajita\" });\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// Change 2: the quoted paper title\n//   \"Context Sensitive Access Control in Smart Home Environments\"\n// (straight quotes, quotes INSIDE the hyperlink run) becomes\n//   \u201cContext Sensitive Access Control in Smart Home Environments\u201d\n// (curly quotes, quotes OUTSIDE the hyperlink).\n// ---------------------------------------------------------------------\n\n// Step 1: strip the leading straight quote from the hyperlink text\n// (partial-run edit so the hyperlink run keeps its Hyperlink style).\nconst leadQuote = body.search('\"Context', { matchCase: true });\nleadQuote.load(\"items\");\nawait context.sync();\nif (leadQuote.items.length > 0) {\n  leadQuote.items[0].insertText(\"Context\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Step 2: strip the trailing straight quote from the hyperlink text.\nconst trailQuote = body.search('Environments\"', { matchCase: true });\ntrailQuote.load(\"items\");\nawait context.sync();\nif (trailQuote.items.length > 0) {\n  trailQuote.items[0].insertText(\"Environments\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Step 3: add the curly opening quote right after \"published paper: \"\n// (merges into the plain-text run that precedes the hyperlink).\nconst beforeLink = body.search(\"published paper: \", { matchCase: true });\nbeforeLink.load(\"items\");\nawait context.sync();\nif (beforeLink.items.length > 0) {\n  beforeLink.items[0].insertText(\"\\u201c\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Step 4: add the curly closing quote right after the hyperlink text\n// (creates a new run outside the hyperlink, inheriting the hyperlink\n// look via rStyle/underline-none, matching the target markup).\nconst linkText = body.search(\n  \"Context Sensitive Access Control in Smart Home Environments\",\n  { matchCase: true }\n);\nlinkText.load(\"items\");\nawait context.sync();\nif (linkText.items.length > 0) {\n  linkText.items[0].insertText(\"\\u201d\", Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: \"...specific technologies like Airflow and Google Cloud\n# Platform.\" -> \"...specific technologies like Apache Airflow and\n# Google Cloud Platform.\"  (rendered as three runs: \"...like \",\n# \"Apache \", \"Airflow and Google Cloud Platform.\")\n# ---------------------------------------------------------------------\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.MatchCase = $true\n$found = $find.Execute(\"Airflow and Google Cloud Platform.\")\nif ($found) {\n    $airflowStart = $range.Start\n    $range.InsertBefore(\"Apache \")\n\n    # \"Apache \" now occupies the 7 characters right at $airflowStart;\n    # touch its font so the engine splits it into its own run.\n    $apacheRange = $d.Range($airflowStart, $airflowStart + 7)\n    $apacheRange.Font.Name = \"Aparajita\"\n\n    # The trailing \"Airflow and Google Cloud Platform.\" text follows\n    # immediately; touch its font too so it becomes a separate run.\n    $airflowRange = $d.Range($airflowStart + 7, $airflowStart + 7 + 35)\n    $airflowRange.Font.Name = \"Aparajita\"\n}\n\n# ---------------------------------------------------------------------\n# Change 2: the quoted paper title\n#   \"Context Sensitive Access Control in Smart Home Environments\"\n# (straight quotes, quotes INSIDE the hyperlink run) becomes\n#   \"Context Sensitive Access Control in Smart Home Environments\"\n# (curly quotes, quotes OUTSIDE the hyperlink).\n# ---------------------------------------------------------------------\n\n# Step 1: delete the leading straight quote from the hyperlink text\n# (single-character Delete() so the hyperlink run keeps its formatting).\n$r1 = $d.Content\n$f1 = $r1.Find\n$f1.ClearFormatting()\n$f1.MatchCase = $true\nif ($f1.Execute('\"Context')) {\n    $d.Range($r1.Start, $r1.Start + 1).Delete()\n}\n\n# Step 2: delete the trailing straight quote from the hyperlink text.\n$r2 = $d.Content\n$f2 = $r2.Find\n$f2.ClearFormatting()\n$f2.MatchCase = $true\nif ($f2.Execute('Environments\"')) {\n    $d.Range($r2.End - 1, $r2.End).Delete()\n}\n\n# Step 3: insert the curly opening quote right after \"published paper: \"\n# (merges into the plain-text run that precedes the hyperlink).\n$r3 = $d.Content\n$f3 = $r3.Find\n$f3.ClearFormatting()\n$f3.MatchCase = $true\nif ($f3.Execute(\"published paper: \")) {\n    $r3.InsertAfter([char]0x201C)\n}\n\n# Step 4: insert the curly closing quote right after the hyperlink text\n# (creates a new run outside the hyperlink, inheriting the hyperlink\n# look via rStyle/underline-none, matching the target markup).\n$r4 = $d.Content\n$f4 = $r4.Find\n$f4.ClearFormatting()\n$f4.MatchCase = $true\nif ($f4.Execute(\"Context Sensitive Access Control in Smart Home Environments\")) {\n    $r4.InsertAfter([char]0x201D)\n}\n"}
